$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 14).Value = 1.32
$ws.Cells.Item(2, 16).Value = 1.32
$ws.Cells.Item(2, 18).Value = 1.21
$ws.Cells.Item(3, 14).Value = 1.29
$ws.Cells.Item(3, 16).Value = 1.28
$ws.Cells.Item(4, 17).Value = 1.76
$ws.Cells.Item(5, 12).Value = 1.18
$ws.Cells.Item(6, 6).Value = 1.99
$ws.Cells.Item(6, 7).Value = 2.2
$ws.Cells.Item(6, 8).Value = 3.5
$ws.Cells.Item(6, 9).Value = 5.1
$ws.Cells.Item(6, 10).Value = 3.25
$ws.Cells.Item(6, 11).Value = 3.85
$ws.Cells.Item(6, 12).Value = 1.43
$ws.Cells.Item(6, 13).Value = 1.07
$ws.Cells.Item(6, 14).Value = 2.94
$ws.Cells.Item(6, 15).Value = 1.35
$ws.Cells.Item(6, 16).Value = 1.77
$ws.Cells.Item(6, 17).Value = 2.02
$ws.Cells.Item(6, 18).Value = 1.29
$ws.Cells.Item(6, 19).Value = 3
$ws.Cells.Item(6, 20).Value = 1.83
$ws.Cells.Item(6, 21).Value = 1.96
$ws.Cells.Item(6, 22).Value = 1.24
$ws.Cells.Item(6, 23).Value = 1.83
$ws.Cells.Item(6, 24).Value = 13.5
$ws.Cells.Item(6, 25).Value = 14.5
$ws.Cells.Item(6, 26).Value = 980
$ws.Cells.Item(6, 27).Value = 110
$ws.Cells.Item(6, 28).Value = 9
$ws.Cells.Item(6, 29).Value = 8.4
$ws.Cells.Item(6, 30).Value = 18
$ws.Cells.Item(6, 31).Value = 60
$ws.Cells.Item(6, 32).Value = 13
$ws.Cells.Item(6, 33).Value = 11.5
$ws.Cells.Item(6, 34).Value = 21
$ws.Cells.Item(6, 35).Value = 75
$ws.Cells.Item(6, 36).Value = 27
$ws.Cells.Item(6, 37).Value = 25
$ws.Cells.Item(6, 38).Value = 55
$ws.Cells.Item(6, 39).Value = 140
$ws.Cells.Item(6, 40).Value = 18.5
$ws.Cells.Item(6, 41).Value = 80
$ws.Cells.Item(7, 12).Value = 1.46
$ws.Cells.Item(7, 13).Value = 1.09
$ws.Cells.Item(7, 14).Value = 3.1
$ws.Cells.Item(7, 15).Value = 1.4
$ws.Cells.Item(7, 18).Value = 1.27
$ws.Cells.Item(7, 19).Value = 4.1
$ws.Cells.Item(7, 20).Value = 1.86
$ws.Cells.Item(7, 21).Value = 2
$ws.Cells.Item(7, 22).Value = 1.46
$ws.Cells.Item(7, 24).Value = 14
$ws.Cells.Item(7, 25).Value = 12.5
$ws.Cells.Item(7, 26).Value = 24
$ws.Cells.Item(7, 27).Value = 55
$ws.Cells.Item(7, 28).Value = 980
$ws.Cells.Item(7, 29).Value = 7.6
$ws.Cells.Item(7, 30).Value = 980
$ws.Cells.Item(7, 31).Value = 980
$ws.Cells.Item(7, 32).Value = 980
$ws.Cells.Item(7, 33).Value = 980
$ws.Cells.Item(7, 34).Value = 23
$ws.Cells.Item(7, 35).Value = 65
$ws.Cells.Item(7, 36).Value = 980
$ws.Cells.Item(7, 37).Value = 980
$ws.Cells.Item(7, 38).Value = 50
$ws.Cells.Item(7, 39).Value = 140
$ws.Cells.Item(7, 40).Value = 40
$ws.Cells.Item(7, 41).Value = 980
$ws.Cells.Item(8, 12).Value = 1.68
$ws.Cells.Item(8, 13).Value = 1.13
$ws.Cells.Item(8, 14).Value = 2.26
$ws.Cells.Item(8, 15).Value = 1.68
$ws.Cells.Item(8, 18).Value = 1.14
$ws.Cells.Item(8, 19).Value = 6.8
$ws.Cells.Item(8, 20).Value = 2.32
$ws.Cells.Item(8, 21).Value = 1.65
$ws.Cells.Item(8, 22).Value = 1.6
$ws.Cells.Item(8, 23).Value = 1.35
$ws.Cells.Item(8, 24).Value = 7.2
$ws.Cells.Item(8, 25).Value = 6.8
$ws.Cells.Item(8, 26).Value = 14.5
$ws.Cells.Item(8, 27).Value = 980
$ws.Cells.Item(8, 28).Value = 9.199999999999999
$ws.Cells.Item(8, 29).Value = 7.2
$ws.Cells.Item(8, 30).Value = 980
$ws.Cells.Item(8, 31).Value = 980
$ws.Cells.Item(8, 32).Value = 980
$ws.Cells.Item(8, 33).Value = 22
$ws.Cells.Item(8, 34).Value = 980
$ws.Cells.Item(8, 35).Value = 110
$ws.Cells.Item(8, 36).Value = 110
$ws.Cells.Item(8, 37).Value = 75
$ws.Cells.Item(8, 38).Value = 140
$ws.Cells.Item(8, 39).Value = 330
$ws.Cells.Item(8, 40).Value = 150
$ws.Cells.Item(8, 41).Value = 55
$ws.Cells.Item(9, 10).Value = 1.12
$ws.Cells.Item(9, 12).Value = 1.01
$ws.Cells.Item(9, 13).Value = 1.01
$ws.Cells.Item(9, 14).Value = 1.89
$ws.Cells.Item(9, 15).Value = 1.01
$ws.Cells.Item(9, 16).Value = 1.3
$ws.Cells.Item(9, 17).Value = 1.02
$ws.Cells.Item(9, 18).Value = 1.11
$ws.Cells.Item(9, 19).Value = 1.01
$ws.Cells.Item(9, 20).Value = 1.01
$ws.Cells.Item(9, 21).Value = 1.01
$ws.Cells.Item(9, 22).Value = 1.01
$ws.Cells.Item(9, 23).Value = 1.01
$ws.Cells.Item(9, 24).Value = 1000
$ws.Cells.Item(9, 25).Value = 1000
$ws.Cells.Item(9, 26).Value = 1000
$ws.Cells.Item(9, 27).Value = 1000
$ws.Cells.Item(9, 28).Value = 1000
$ws.Cells.Item(9, 29).Value = 1000
$ws.Cells.Item(9, 30).Value = 1000
$ws.Cells.Item(9, 31).Value = 1000
$ws.Cells.Item(9, 32).Value = 1000
$ws.Cells.Item(9, 33).Value = 1000
$ws.Cells.Item(9, 34).Value = 1000
$ws.Cells.Item(9, 35).Value = 1000
$ws.Cells.Item(9, 36).Value = 1000
$ws.Cells.Item(9, 37).Value = 1000
$ws.Cells.Item(9, 38).Value = 1000
$ws.Cells.Item(9, 39).Value = 1000
$ws.Cells.Item(9, 40).Value = 1000
$ws.Cells.Item(9, 41).Value = 1000
$ws.Cells.Item(10, 6).Value = 2.2
$ws.Cells.Item(10, 7).Value = 2.4
$ws.Cells.Item(10, 8).Value = 3.85
$ws.Cells.Item(10, 9).Value = 4.3
$ws.Cells.Item(10, 11).Value = 3.3
$ws.Cells.Item(10, 12).Value = 1.01
$ws.Cells.Item(10, 13).Value = 1.01
$ws.Cells.Item(10, 14).Value = 1.55
$ws.Cells.Item(10, 15).Value = 1.02
$ws.Cells.Item(10, 16).Value = 1.55
$ws.Cells.Item(10, 17).Value = 2.5
$ws.Cells.Item(10, 18).Value = 1.17
$ws.Cells.Item(10, 19).Value = 4.3
$ws.Cells.Item(10, 20).Value = 1.74
$ws.Cells.Item(10, 21).Value = 1.53
$ws.Cells.Item(10, 22).Value = 1.3
$ws.Cells.Item(10, 23).Value = 1.71
$ws.Cells.Item(10, 24).Value = 13
$ws.Cells.Item(10, 25).Value = 16
$ws.Cells.Item(10, 26).Value = 980
$ws.Cells.Item(10, 27).Value = 1000
$ws.Cells.Item(10, 28).Value = 10.5
$ws.Cells.Item(10, 29).Value = 10.5
$ws.Cells.Item(10, 30).Value = 980
$ws.Cells.Item(10, 31).Value = 100
$ws.Cells.Item(10, 32).Value = 980
$ws.Cells.Item(10, 33).Value = 17.5
$ws.Cells.Item(10, 34).Value = 980
$ws.Cells.Item(10, 35).Value = 1000
$ws.Cells.Item(10, 36).Value = 980
$ws.Cells.Item(10, 37).Value = 980
$ws.Cells.Item(10, 38).Value = 90
$ws.Cells.Item(10, 39).Value = 1000
$ws.Cells.Item(10, 40).Value = 1000
$ws.Cells.Item(10, 41).Value = 1000
$ws.Cells.Item(11, 7).Value = 2.2
$ws.Cells.Item(11, 8).Value = 4
$ws.Cells.Item(11, 12).Value = 1.49
$ws.Cells.Item(11, 13).Value = 1.09
$ws.Cells.Item(11, 14).Value = 3.05
$ws.Cells.Item(11, 15).Value = 1.42
$ws.Cells.Item(11, 18).Value = 1.26
$ws.Cells.Item(11, 19).Value = 4.1
$ws.Cells.Item(11, 20).Value = 1.93
$ws.Cells.Item(11, 21).Value = 1.9
$ws.Cells.Item(11, 22).Value = 1.27
$ws.Cells.Item(11, 23).Value = 1.83
$ws.Cells.Item(11, 24).Value = 12
$ws.Cells.Item(11, 25).Value = 13.5
$ws.Cells.Item(11, 26).Value = 980
$ws.Cells.Item(11, 27).Value = 120
$ws.Cells.Item(11, 28).Value = 8.4
$ws.Cells.Item(11, 29).Value = 8
$ws.Cells.Item(11, 30).Value = 18.5
$ws.Cells.Item(11, 31).Value = 65
$ws.Cells.Item(11, 32).Value = 13
$ws.Cells.Item(11, 33).Value = 11.5
$ws.Cells.Item(11, 34).Value = 22
$ws.Cells.Item(11, 35).Value = 80
$ws.Cells.Item(11, 36).Value = 980
$ws.Cells.Item(11, 37).Value = 27
$ws.Cells.Item(11, 38).Value = 980
$ws.Cells.Item(11, 39).Value = 170
$ws.Cells.Item(11, 40).Value = 22
$ws.Cells.Item(11, 41).Value = 100
$ws.Cells.Item(12, 12).Value = 1.01
$ws.Cells.Item(12, 13).Value = 1.01
$ws.Cells.Item(12, 14).Value = 1.01
$ws.Cells.Item(12, 15).Value = 1.34
$ws.Cells.Item(12, 18).Value = 1.18
$ws.Cells.Item(12, 19).Value = 1.01
$ws.Cells.Item(12, 20).Value = 1.01
$ws.Cells.Item(12, 21).Value = 1.01
$ws.Cells.Item(12, 22).Value = 1.01
$ws.Cells.Item(12, 23).Value = 1.01
$ws.Cells.Item(12, 24).Value = 1000
$ws.Cells.Item(12, 25).Value = 1000
$ws.Cells.Item(12, 26).Value = 1000
$ws.Cells.Item(12, 27).Value = 1000
$ws.Cells.Item(12, 28).Value = 1000
$ws.Cells.Item(12, 29).Value = 1000
$ws.Cells.Item(12, 30).Value = 1000
$ws.Cells.Item(12, 31).Value = 1000
$ws.Cells.Item(12, 32).Value = 1000
$ws.Cells.Item(12, 33).Value = 1000
$ws.Cells.Item(12, 34).Value = 1000
$ws.Cells.Item(12, 35).Value = 1000
$ws.Cells.Item(12, 36).Value = 1000
$ws.Cells.Item(12, 37).Value = 1000
$ws.Cells.Item(12, 38).Value = 1000
$ws.Cells.Item(12, 39).Value = 1000
$ws.Cells.Item(12, 40).Value = 1000
$ws.Cells.Item(12, 41).Value = 1000